# Insert a new daily price record as the newest row (row 8), pushing all
# the existing rows from 8..78 down to 9..79. The shifted rows keep their
# original values automatically via the Insert, so we only need to
# populate the freshly inserted row 8 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("8:8").Insert()

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44552
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100102
$ws.Range("H8").Value = "Cítricos"
$ws.Range("I8").Value = 100102005
$ws.Range("J8").Value = "Naranja"
$ws.Range("K8").Value = "Valencia"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 700
$ws.Range("O8").Value = 750
$ws.Range("P8").Value = 725
$ws.Range("Q8").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R8").Value = "Región de Coquimbo"
$ws.Range("S8").Value = 725
$ws.Range("T8").Value = 1
